$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy formatting (bold, centered, boxed) from the existing
# header cell H1, then set the new header text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J (rows 2-13)
$iValues = @(1,1,1,1,1,1,1,1,1,6,7,5)
$jValues = @(4,5,5,5,4,6,4,3,5,7,7,5)

for ($idx = 0; $idx -lt 12; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
